# Apply changes to sheet "保險" (insurance, sheet index 6) and
# "債權" (claim, sheet index 7): add the standard trailing metadata
# columns (property_category/category, date, legislator_name,
# legislator_id, source_file, index) that the other sheets already have,
# and turn the header row into proper field-name labels.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "保險" (insurance) -> Worksheets index 6
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Header row: field names
$ws6.Range("B1").Value = "company"
$ws6.Range("C1").Value = "name"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "property_category"
$ws6.Range("F1").Value = "category"
$ws6.Range("G1").Value = "date"
$ws6.Range("H1").Value = "legislator_name"
$ws6.Range("I1").Value = "legislator_id"
$ws6.Range("J1").Value = "source_file"
$ws6.Range("K1").Value = "index"

# match header style (bold, centered, bordered) for the newly added cells
$hdrNew6 = $ws6.Range("F1:K1")
$hdrNew6.Font.Bold = $true
$hdrNew6.HorizontalAlignment = -4108
$hdrNew6.VerticalAlignment = -4160
$hdrNew6.Borders.LineStyle = 1

# Row 2 (index 135)
$ws6.Range("E2").Value = "insurance"
$ws6.Range("F2").Value = "normal"
$ws6.Range("G2").NumberFormat = "@"
$ws6.Range("G2").Value = "2011-11-21"
$ws6.Range("G2").Style = "Normal"
$ws6.Range("H2").Value = "羅淑蕾"
$ws6.Range("I2").Value = 1638
$ws6.Range("J2").Value = "tmpe5cc1"
$ws6.Range("K2").Value = 135

# Row 3 (index 136)
$ws6.Range("E3").Value = "insurance"
$ws6.Range("F3").Value = "normal"
$ws6.Range("G3").NumberFormat = "@"
$ws6.Range("G3").Value = "2011-11-21"
$ws6.Range("G3").Style = "Normal"
$ws6.Range("H3").Value = "羅淑蕾"
$ws6.Range("I3").Value = 1638
$ws6.Range("J3").Value = "tmpe5cc1"
$ws6.Range("K3").Value = 136

# ---------------------------------------------------------------
# Sheet "債權" (claim) -> Worksheets index 7
# ---------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Header row: field names
$ws7.Range("B1").Value = "species"
$ws7.Range("C1").Value = "owner"
$ws7.Range("D1").Value = "debtor"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"

# match header style (bold, centered, bordered) for the newly added cells
$hdrNew7 = $ws7.Range("H1:N1")
$hdrNew7.Font.Bold = $true
$hdrNew7.HorizontalAlignment = -4108
$hdrNew7.VerticalAlignment = -4160
$hdrNew7.Borders.LineStyle = 1

# Row 2 (index 141)
$ws7.Range("B2").Value = "設定抵押"
$ws7.Range("D2").Value = "王國道臺北市文山區興隆路"
$ws7.Range("F2").Value = "93年12月01日"
$ws7.Range("G2").Value = "借款"
$ws7.Range("H2").Value = "claim"
$ws7.Range("I2").Value = "normal"
$ws7.Range("J2").NumberFormat = "@"
$ws7.Range("J2").Value = "2011-11-21"
$ws7.Range("J2").Style = "Normal"
$ws7.Range("K2").Value = "羅淑蕾"
$ws7.Range("L2").Value = 1638
$ws7.Range("M2").Value = "tmpe5cc1"
$ws7.Range("N2").Value = 141

# Row 3 (index 142)
$ws7.Range("B3").Value = "—般借款"
$ws7.Range("D3").Value = "啟富建設有限公司臺北市中山區建國北路2段92號9樓"
$ws7.Range("F3").Value = "96年05月01日"
$ws7.Range("G3").Value = "借款"
$ws7.Range("H3").Value = "claim"
$ws7.Range("I3").Value = "normal"
$ws7.Range("J3").NumberFormat = "@"
$ws7.Range("J3").Value = "2011-11-21"
$ws7.Range("J3").Style = "Normal"
$ws7.Range("K3").Value = "羅淑蕾"
$ws7.Range("L3").Value = 1638
$ws7.Range("M3").Value = "tmpe5cc1"
$ws7.Range("N3").Value = 142
